$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -0.3

$ws.Range("B3").Value = 2
$ws.Range("C3").Value = -0.7

$ws.Range("C4").Value = 20.9

$ws.Range("C7").Value = 8.4

$ws.Range("C10").Value = -0.4

$ws.Range("C11").Value = 2.1

$ws.Range("C13").Value = -8

$ws.Range("C20").Value = -1.1

$ws.Range("C22").Value = -9.199999999999999
